# "INFORMATIVO 16 GRE" workbook update
# The plate "KIQ-8720" (rows 30:35, column A) was actually "KZQ-8720".
# Fix the typo in those six cells; Excel drops the cells' existing
# (bordered/centered) style when the content is corrected, so we clear
# the formatting on that range before writing the corrected value.
# Also move the current on-screen selection to F40 (was G65).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$plateRange = $ws.Range("A30:A35")
$plateRange.ClearFormats()

$ws.Range("A30").Value = "KZQ-8720"
$ws.Range("A31").Value = "KZQ-8720"
$ws.Range("A32").Value = "KZQ-8720"
$ws.Range("A33").Value = "KZQ-8720"
$ws.Range("A34").Value = "KZQ-8720"
$ws.Range("A35").Value = "KZQ-8720"

$ws.Range("F40").Select()
